$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data rows for sub_004 and sub_005
$ws.Range("A5").Value = "sub_004"
$ws.Range("B5").Value = $false

$ws.Range("A6").Value = "sub_005"
$ws.Range("B6").Value = $false

# Update the saved selection/active cell display info
$ws.Range("C12").Select()
